$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text-formatted numbers (e.g. "63.15", "213.38").
# Force it to Text format so Excel does not silently convert these
# look-alike numeric strings into real numbers (which would lose trailing
# zeros / change precision), matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.684.90"
$ws.Range("E2").Value = "  +1.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.635.42"
$ws.Range("E3").Value = "  +1.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.38"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.09%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.98%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.58%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.63%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.02"
$ws.Range("E10").Value = "  +2.74%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  +2.58%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.863.05"
$ws.Range("E12").Value = "  +1.41%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.633.52"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.04%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  +1.61%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.655.29"
$ws.Range("E16").Value = "  +1.25%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.15"
$ws.Range("E17").Value = "  +2.03%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  +0.44%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "208.38"

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.54%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  +0.44%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +0.97%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - Monero
$ws.Range("D25").Value = "145.92"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.65%  "

# Rows 28 and 29 swap: Cosmos <-> EthereumClassic (with updated values)
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.36"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "6.67"
$ws.Range("E29").Value = "  +1.46%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +5.66%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.27%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.90%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.48%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.06%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.42%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.163.03"
$ws.Range("E36").Value = "  +0.04%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +1.14%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +1.81%  "

# Row 39 - PaxDollar
$ws.Range("E39").Value = "  +0.08%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  -0.20%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.796"
$ws.Range("E42").Value = "  +1.08%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +2.68%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.772.21"
$ws.Range("E44").Value = "  +1.18%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.40"
$ws.Range("E45").Value = "  +0.64%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +0.26%  "

# Row 47 - Aave
$ws.Range("E47").Value = "  +0.35%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +7.92%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  +0.70%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "7.59"
$ws.Range("E50").Value = "  +4.67%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "0.411"
$ws.Range("E51").Value = "  +1.05%  "
